$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update cell values to match the reshuffled name list ---
$ws.Range('B1').Value = '溫銘'
$ws.Range('D1').Value = '葉雪梅'
$ws.Range('E1').Value = '張志謀'
$ws.Range('B2').Value = '何惠珍'
$ws.Range('D2').Value = '馮偉初'
$ws.Range('E2').Value = '張志謀'
$ws.Range('B3').Value = '溫泉德'
$ws.Range('D3').Value = '葉雪容'
$ws.Range('E3').Value = '*Yellow'
$ws.Range('B4').Value = '溫子茹'
$ws.Range('D4').Value = '葉紅志'
$ws.Range('E4').Value = 'Ms Tong'
$ws.Range('B5').Value = '溫國郷'
$ws.Range('D5').Value = '廖練雲'
$ws.Range('E5').Value = 'Denis'
$ws.Range('G5').Value = '*Anthony '
$ws.Range('B6').Value = '何祥增'
$ws.Range('D6').Value = '廖苑雲'
$ws.Range('E6').Value = 'Begger'
$ws.Range('G6').Value = 'Tim'
$ws.Range('B7').Value = '溫素芳'
$ws.Range('D7').Value = '廖江真'
$ws.Range('E7').Value = 'Sol Bread'
$ws.Range('G7').Value = 'Patrick'
$ws.Range('B8').Value = '廖日雲'
$ws.Range('C8').Value = '廖鋼基'
$ws.Range('D8').Value = '廖睦堯'
$ws.Range('E8').Value = 'Fai Chi'
$ws.Range('G8').Value = 'See Fu'
$ws.Range('B9').Value = '廖雲基'
$ws.Range('C9').Value = '廖偉良'
$ws.Range('D9').Value = '廖維華'
$ws.Range('E9').Value = '#Somingtat'
$ws.Range('G9').Value = 'Yvoone'
$ws.Range('B10').Value = '廖彩雲'
$ws.Range('C10').Value = '廖素琼'
$ws.Range('D10').Value = '廖維華'
$ws.Range('E10').Value = 'SomingtatW'
$ws.Range('G10').Value = 'Pui'
$ws.Range('B11').Value = '廖志'
$ws.Range('C11').Value = '黃錦萍'
$ws.Range('D11').Value = '廖富盛'
$ws.Range('E11').Value = 'Jasper'
$ws.Range('G11').Value = 'Ocean'
$ws.Range('B12').Value = '黃塋塋'
$ws.Range('C12').Value = 'AR NE'
$ws.Range('D12').Value = '廖富盛'
$ws.Range('E12').Value = 'Tin Shing'
$ws.Range('G12').Value = '**Lam Kei'
$ws.Range('D13').Value = '廖富盛'
$ws.Range('G18').Value = 'Wong Ming'
$ws.Range('G19').Value = 'Wong Kei'
$ws.Range('G20').Value = 'Zuey Tsui'
$ws.Range('G21').Value = 'Guanglei'
$ws.Range('G22').Value = 'Bean Man'
$ws.Range('G23').Value = 'Foo Kwai'
$ws.Range('G24').Value = 'Ellen '

# --- Clear cells that no longer hold data ---
$ws.Range("G13").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("E21").ClearContents()

# --- Shaded (gray125) highlight follows specific entries ---
$ws.Range("E4").Interior.Pattern = -4142
$ws.Range("E10").Interior.Pattern = -4142
$ws.Range("E3").Interior.Pattern = 17
$ws.Range("G5").Interior.Pattern = 17

# --- Update the active selection ---
$ws.Range("E13").Select()
